# Commit message: "remove obsolete columns from data"
#
# The "Brands" worksheet contained columns for brand, shop, type, quality,
# origin, package, year, kilogram, price, team. The origin, package and
# team columns are obsolete and are removed entirely (shifting the
# columns to their right one step to the left). This also causes the
# now-unused shared strings (origin, package, Spain, Net, team, unknown)
# to disappear from the workbook's shared string table, and the "Fruits"
# worksheet (which never referenced those strings) ends up pointing at
# the new, lower string indices after the table is compacted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brands")

# Delete right-to-left so remaining column letters keep their meaning
# while we work through the list.
$null = $ws.Range("J:J").EntireColumn.Delete()   # team
$null = $ws.Range("F:F").EntireColumn.Delete()   # package
$null = $ws.Range("E:E").EntireColumn.Delete()   # origin

# Restore a sensible selection on the sheet (mirrors the state the
# workbook was left in after the edit).
$null = $ws.Range("H13").Select()
